# Fix add record to add multiple records in session
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Update existing row 2 values
$ws.Range("G2").Value = 7107
$ws.Range("K2").Value = "wells fargo"

# Update existing row 3 values
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 35803.4
$ws.Range("L3").Value = "Wages, tips, other comp."

# Add new row 4
$ws.Range("A4").Value = "16afe2b7-b1ee-4c6a-ad30-699637f33f6f"
# Enter the date as a text formula then flatten to a literal value so it is
# stored as plain text "01/31/2023" (not auto-converted to a date serial)
# and without leaving a lingering text-number-format style override.
$ws.Range("B4").Formula = '="01/31/2023"'
$ws.Range("B4").Copy()
$ws.Range("B4").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("H4").Value = "tw22"

# Add new row 5
$ws.Range("A5").Value = "b6818292-a689-4a17-8b0f-ea2f05de255a"
$ws.Range("B5").Formula = '="01/31/2023"'
$ws.Range("B5").Copy()
$ws.Range("B5").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("H5").Value = "tw23"

# Apply style from A1 (header bold/border style) to new ID cells A4:A5 to match style s="1"
$ws.Range("A1").Copy()
$ws.Range("A4:A5").PasteSpecial(-4122)  # xlPasteFormats
